$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (preserve exact formatting incl. trailing zeros / % sign)
# for the numeric-looking Price/Volume columns before assigning their new values.
$numericCells = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "D19", "E19", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "E25", "D26", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "E44", "D45", "E45", "D46", "E46", "E47", "E48", "D49", "E49", "E50", "E51")
foreach ($addr in $numericCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Updated price / volume(1h) figures
$ws.Range("D2").Value = "325.67"
$ws.Range("E2").Value = "-1.11%"
$ws.Range("D3").Value = "39.53"
$ws.Range("E3").Value = "-1.10%"
$ws.Range("D4").Value = "5.660"
$ws.Range("E4").Value = "6.20%"
$ws.Range("D5").Value = "0.08030"
$ws.Range("E5").Value = "-0.91%"
$ws.Range("D6").Value = "2.041"
$ws.Range("E6").Value = "6.12%"
$ws.Range("D7").Value = "8.622"
$ws.Range("E7").Value = "-0.40%"
$ws.Range("D8").Value = "4.476"
$ws.Range("E8").Value = "-0.95%"
$ws.Range("E9").Value = "-1.47%"
$ws.Range("D10").Value = "0.9222"
$ws.Range("E10").Value = "-2.21%"
$ws.Range("D11").Value = "0.1245"
$ws.Range("E11").Value = "-8.35%"
$ws.Range("D12").Value = "0.1965"
$ws.Range("E12").Value = "-0.18%"
$ws.Range("D13").Value = "8.718"
$ws.Range("E13").Value = "21.15%"
$ws.Range("D14").Value = "0.09192"
$ws.Range("E14").Value = "-1.38%"
$ws.Range("D15").Value = "0.03558"
$ws.Range("E15").Value = "-0.29%"
$ws.Range("D16").Value = "0.1050"
$ws.Range("E16").Value = "9.51%"
$ws.Range("D17").Value = "0.001292"
$ws.Range("E17").Value = "-2.43%"
$ws.Range("D18").Value = "0.006119"
$ws.Range("E18").Value = "-4.83%"
$ws.Range("D19").Value = "3.353"
$ws.Range("E19").Value = "-0.37%"
$ws.Range("E20").Value = "-1.14%"
$ws.Range("D21").Value = "0.1371"
$ws.Range("E21").Value = "2.86%"
$ws.Range("D22").Value = "0.2504"
$ws.Range("E22").Value = "-2.27%"
$ws.Range("D23").Value = "0.04386"
$ws.Range("E23").Value = "-0.91%"
$ws.Range("D24").Value = "0.001260"
$ws.Range("E24").Value = "3.08%"
$ws.Range("E25").Value = "6.50%"
$ws.Range("D26").Value = "0.0001230"
$ws.Range("D39").Value = "0.02515"
$ws.Range("E39").Value = "0.94%"
$ws.Range("D40").Value = "0.05336"
$ws.Range("E40").Value = "2.21%"
$ws.Range("D41").Value = "0.007491"
$ws.Range("E41").Value = "-1.96%"
$ws.Range("D42").Value = "0.009908"
$ws.Range("E42").Value = "9.24%"
$ws.Range("D43").Value = "0.1407"
$ws.Range("E43").Value = "-1.43%"
$ws.Range("E44").Value = "-2.11%"
$ws.Range("D45").Value = "0.01109"
$ws.Range("E45").Value = "2.06%"
$ws.Range("D46").Value = "0.00006691"
$ws.Range("E46").Value = "0.90%"
$ws.Range("E47").Value = "-0.04%"
$ws.Range("E48").Value = "-5.08%"
$ws.Range("D49").Value = "0.002975"
$ws.Range("E49").Value = "-11.06%"
$ws.Range("E50").Value = "-0.04%"
$ws.Range("E51").Value = "-0.04%"

# Row 7/8 coin swap: KuCoinToken now ranked above GateToken
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
